# Applies the updated Price (D) / Volume(1h) (E) figures from the latest
# cryptos data refresh. Values are written as literal text (not numbers),
# matching how the source sheet stores them (e.g. '1.003', '25.722.79').
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    # Leading apostrophe forces Excel to store the content as literal text,
    # exactly like typing '1.003 into a cell -- prevents numeric/date coercion
    # of values such as '10.00', '1.003' or '0.000007990'.
    $ws.Range($Cell).Value = "'" + $Text
}

$ws.Range('D2').Value = '25.722.79'
$ws.Range('E2').Value = '  -5.55%  '
$ws.Range('D3').Value = '1.814.42'
$ws.Range('E3').Value = '  -4.61%  '
Set-TextValue 'D4' '1.003'
$ws.Range('E4').Value = '  +0.08%  '
Set-TextValue 'D5' '278.14'
$ws.Range('E5').Value = '  -9.19%  '
Set-TextValue 'D6' '1.002'
$ws.Range('E6').Value = '  +0.13%  '
Set-TextValue 'D7' '0.4945'
$ws.Range('E7').Value = '  -8.01%  '
Set-TextValue 'D8' '0.3497'
$ws.Range('E8').Value = '  -8.25%  '
Set-TextValue 'D9' '44.29'
$ws.Range('E9').Value = '  -3.73%  '
Set-TextValue 'D10' '0.06612'
$ws.Range('E10').Value = '  -9.31%  '
Set-TextValue 'D11' '20.08'
$ws.Range('E11').Value = '  -9.47%  '
Set-TextValue 'D12' '0.8458'
$ws.Range('E12').Value = '  -6.30%  '
Set-TextValue 'D13' '0.07813'
$ws.Range('E13').Value = '  -4.85%  '
$ws.Range('D14').Value = '1.810.75'
$ws.Range('E14').Value = '  +64.42%  '
Set-TextValue 'D15' '5.033'
$ws.Range('E15').Value = '  -5.69%  '
Set-TextValue 'D16' '87.35'
$ws.Range('E16').Value = '  -9.19%  '
Set-TextValue 'D17' '1.003'
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('E18').Value = '  -6.10%  '
Set-TextValue 'D19' '1.003'
$ws.Range('E19').Value = '  +0.13%  '
Set-TextValue 'D20' '0.000007990'
$ws.Range('E20').Value = '  -7.59%  '
$ws.Range('D21').Value = '25.783.83'
$ws.Range('E21').Value = '  -5.44%  '
$ws.Range('E22').Value = '  -5.75%  '
Set-TextValue 'D23' '10.00'
$ws.Range('E23').Value = '  -7.18%  '
Set-TextValue 'D24' '6.087'
$ws.Range('E24').Value = '  -6.29%  '
Set-TextValue 'D25' '141.26'
$ws.Range('E25').Value = '  -5.57%  '
Set-TextValue 'D26' '2.131'
$ws.Range('E26').Value = '  -6.92%  '
Set-TextValue 'D27' '1.663'
$ws.Range('E27').Value = '  -4.64%  '
Set-TextValue 'D28' '16.81'
$ws.Range('E28').Value = '  -8.34%  '
Set-TextValue 'D29' '108.77'
$ws.Range('E29').Value = '  -6.78%  '
Set-TextValue 'D30' '4.295'
$ws.Range('E30').Value = '  -10.60%  '
Set-TextValue 'D31' '4.216'
$ws.Range('E31').Value = '  -11.52%  '
Set-TextValue 'D32' '0.08763'
$ws.Range('E32').Value = '  -4.97%  '
Set-TextValue 'D33' '0.04794'
$ws.Range('E33').Value = '  -5.41%  '
Set-TextValue 'D34' '0.7405'
$ws.Range('E34').Value = '  -11.16%  '
Set-TextValue 'D35' '2.882'
$ws.Range('E35').Value = '  -3.70%  '
Set-TextValue 'D36' '1.129'
$ws.Range('E36').Value = '  -7.33%  '
Set-TextValue 'D37' '1.003'
$ws.Range('E37').Value = '  +0.32%  '
Set-TextValue 'D38' '3.060'
$ws.Range('E38').Value = '  -8.55%  '
Set-TextValue 'D39' '2.454'
$ws.Range('E39').Value = '  -8.26%  '
Set-TextValue 'D40' '0.5310'
$ws.Range('E40').Value = '  -8.05%  '
Set-TextValue 'D41' '0.01849'
$ws.Range('E41').Value = '  -7.85%  '
Set-TextValue 'D42' '0.9683'
$ws.Range('E42').Value = '  -9.97%  '
Set-TextValue 'D43' '113.69'
$ws.Range('E43').Value = '  -2.77%  '
Set-TextValue 'D44' '6.199'
$ws.Range('E44').Value = '  -6.11%  '
Set-TextValue 'D45' '8.158'
$ws.Range('E45').Value = '  -12.93%  '
Set-TextValue 'D46' '0.4704'
$ws.Range('E46').Value = '  -5.14%  '
Set-TextValue 'D48' '0.1381'
$ws.Range('E48').Value = '  -9.32%  '
Set-TextValue 'D49' '9.171'
$ws.Range('E49').Value = '  -9.18%  '
Set-TextValue 'D50' '35.84'
$ws.Range('E50').Value = '  -6.63%  '
Set-TextValue 'D51' '0.05877'
$ws.Range('E51').Value = '  -4.81%  '
